$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 16,51

$arr[0,0] = 111814212
$arr[0,1] = 89405
$arr[0,2] = "Ovaliderad"
$arr[0,3] = "NT"
$arr[0,4] = 1202
$arr[0,5] = "Ullticka"
$arr[0,6] = "Phellinidium ferrugineofuscum"
$arr[0,7] = "(P.Karst.) Fiasson & Niemelä"
$arr[0,8] = $null
$arr[0,9] = $null
$arr[0,10] = $null
$arr[0,11] = $null
$arr[0,12] = $null
$arr[0,13] = $null
$arr[0,14] = $null
$arr[0,15] = "Matsdal, granskog, Ås lm"
$arr[0,16] = 540635.9369002836
$arr[0,17] = 7247595.565451766
$arr[0,18] = 10
$arr[0,19] = "Västerbotten"
$arr[0,20] = "Vilhelmina"
$arr[0,21] = "Åsele lappmark"
$arr[0,22] = "Vilhelmina"
$arr[0,23] = $null
$arr[0,24] = "2023-08-18"
$arr[0,25] = "11:00"
$arr[0,26] = "2023-08-18"
$arr[0,27] = "11:00"
$arr[0,28] = $null
$arr[0,29] = $false
$arr[0,30] = $false
$arr[0,31] = $null
$arr[0,32] = $false
$arr[0,33] = $null
$arr[0,34] = $null
$arr[0,35] = $null
$arr[0,36] = $null
$arr[0,37] = $null
$arr[0,38] = $null
$arr[0,39] = $null
$arr[0,40] = $null
$arr[0,41] = $null
$arr[0,42] = $null
$arr[0,43] = $null
$arr[0,44] = $null
$arr[0,45] = $null
$arr[0,46] = $null
$arr[0,47] = $null
$arr[0,48] = "Roger Olofsson"
$arr[0,49] = "Roger Olofsson"
$arr[0,50] = $null

$arr[1,0] = 111814047
$arr[1,1] = 90087
$arr[1,2] = "Ovaliderad"
$arr[1,3] = "LC"
$arr[1,4] = 3298
$arr[1,5] = "Trådticka"
$arr[1,6] = "Climacocystis borealis"
$arr[1,7] = "(Fr.) Kotl. & Pouzar"
$arr[1,8] = $null
$arr[1,9] = $null
$arr[1,10] = $null
$arr[1,11] = $null
$arr[1,12] = $null
$arr[1,13] = $null
$arr[1,14] = $null
$arr[1,15] = "Matsdal, granskog, Ås lm"
$arr[1,16] = 540633.6855369165
$arr[1,17] = 7247516.598344535
$arr[1,18] = 10
$arr[1,19] = "Västerbotten"
$arr[1,20] = "Vilhelmina"
$arr[1,21] = "Åsele lappmark"
$arr[1,22] = "Vilhelmina"
$arr[1,23] = $null
$arr[1,24] = "2023-08-18"
$arr[1,25] = "11:00"
$arr[1,26] = "2023-08-18"
$arr[1,27] = "11:00"
$arr[1,28] = "rikligt"
$arr[1,29] = $false
$arr[1,30] = $false
$arr[1,31] = $null
$arr[1,32] = $false
$arr[1,33] = $null
$arr[1,34] = $null
$arr[1,35] = $null
$arr[1,36] = $null
$arr[1,37] = $null
$arr[1,38] = $null
$arr[1,39] = $null
$arr[1,40] = $null
$arr[1,41] = $null
$arr[1,42] = $null
$arr[1,43] = $null
$arr[1,44] = $null
$arr[1,45] = $null
$arr[1,46] = $null
$arr[1,47] = $null
$arr[1,48] = "Roger Olofsson"
$arr[1,49] = "Roger Olofsson"
$arr[1,50] = $null

$arr[2,0] = 111813745
$arr[2,1] = 56398
$arr[2,2] = "Ovaliderad"
$arr[2,3] = "NT"
$arr[2,4] = 100109
$arr[2,5] = "Tretåig hackspett"
$arr[2,6] = "Picoides tridactylus"
$arr[2,7] = "(Linnaeus, 1758)"
$arr[2,8] = $null
$arr[2,9] = $null
$arr[2,10] = $null
$arr[2,11] = $null
$arr[2,12] = "färska spår"
$arr[2,13] = $null
$arr[2,14] = $null
$arr[2,15] = "Matsdal, granskog, Ås lm"
$arr[2,16] = 540568.950047517
$arr[2,17] = 7247601.73830481
$arr[2,18] = 10
$arr[2,19] = "Västerbotten"
$arr[2,20] = "Vilhelmina"
$arr[2,21] = "Åsele lappmark"
$arr[2,22] = "Vilhelmina"
$arr[2,23] = $null
$arr[2,24] = "2023-08-18"
$arr[2,25] = "11:00"
$arr[2,26] = "2023-08-18"
$arr[2,27] = "11:00"
$arr[2,28] = $null
$arr[2,29] = $false
$arr[2,30] = $false
$arr[2,31] = $null
$arr[2,32] = $false
$arr[2,33] = $null
$arr[2,34] = $null
$arr[2,35] = "gran"
$arr[2,36] = "Picea abies"
$arr[2,37] = $null
$arr[2,38] = $null
$arr[2,39] = $null
$arr[2,40] = "Picea abies"
$arr[2,41] = $null
$arr[2,42] = $null
$arr[2,43] = $null
$arr[2,44] = $null
$arr[2,45] = $null
$arr[2,46] = $null
$arr[2,47] = $null
$arr[2,48] = "Roger Olofsson"
$arr[2,49] = "Roger Olofsson"
$arr[2,50] = $null

$arr[3,0] = 111814303
$arr[3,1] = 90087
$arr[3,2] = "Ovaliderad"
$arr[3,3] = "LC"
$arr[3,4] = 3298
$arr[3,5] = "Trådticka"
$arr[3,6] = "Climacocystis borealis"
$arr[3,7] = "(Fr.) Kotl. & Pouzar"
$arr[3,8] = $null
$arr[3,9] = $null
$arr[3,10] = $null
$arr[3,11] = $null
$arr[3,12] = $null
$arr[3,13] = $null
$arr[3,14] = $null
$arr[3,15] = "Matsdal, granskog, Ås lm"
$arr[3,16] = 540600.641023421
$arr[3,17] = 7247517.393825463
$arr[3,18] = 10
$arr[3,19] = "Västerbotten"
$arr[3,20] = "Vilhelmina"
$arr[3,21] = "Åsele lappmark"
$arr[3,22] = "Vilhelmina"
$arr[3,23] = $null
$arr[3,24] = "2023-08-18"
$arr[3,25] = "11:00"
$arr[3,26] = "2023-08-18"
$arr[3,27] = "11:00"
$arr[3,28] = $null
$arr[3,29] = $false
$arr[3,30] = $false
$arr[3,31] = $null
$arr[3,32] = $false
$arr[3,33] = $null
$arr[3,34] = $null
$arr[3,35] = $null
$arr[3,36] = $null
$arr[3,37] = $null
$arr[3,38] = $null
$arr[3,39] = $null
$arr[3,40] = $null
$arr[3,41] = $null
$arr[3,42] = $null
$arr[3,43] = $null
$arr[3,44] = $null
$arr[3,45] = $null
$arr[3,46] = $null
$arr[3,47] = $null
$arr[3,48] = "Roger Olofsson"
$arr[3,49] = "Roger Olofsson"
$arr[3,50] = $null

$arr[4,0] = 111814152
$arr[4,1] = 89423
$arr[4,2] = "Ovaliderad"
$arr[4,3] = "NT"
$arr[4,4] = 5432
$arr[4,5] = "Granticka"
$arr[4,6] = "Porodaedalea chrysoloma"
$arr[4,7] = "(Fr.) Fiasson & Niemelä"
$arr[4,8] = $null
$arr[4,9] = $null
$arr[4,10] = $null
$arr[4,11] = $null
$arr[4,12] = $null
$arr[4,13] = $null
$arr[4,14] = $null
$arr[4,15] = "Matsdal, granskog, Ås lm"
$arr[4,16] = 540661.0419420782
$arr[4,17] = 7247564.172119373
$arr[4,18] = 10
$arr[4,19] = "Västerbotten"
$arr[4,20] = "Vilhelmina"
$arr[4,21] = "Åsele lappmark"
$arr[4,22] = "Vilhelmina"
$arr[4,23] = $null
$arr[4,24] = "2023-08-18"
$arr[4,25] = "11:00"
$arr[4,26] = "2023-08-18"
$arr[4,27] = "11:00"
$arr[4,28] = $null
$arr[4,29] = $false
$arr[4,30] = $false
$arr[4,31] = $null
$arr[4,32] = $false
$arr[4,33] = $null
$arr[4,34] = $null
$arr[4,35] = $null
$arr[4,36] = $null
$arr[4,37] = $null
$arr[4,38] = $null
$arr[4,39] = $null
$arr[4,40] = $null
$arr[4,41] = $null
$arr[4,42] = $null
$arr[4,43] = $null
$arr[4,44] = $null
$arr[4,45] = $null
$arr[4,46] = $null
$arr[4,47] = $null
$arr[4,48] = "Roger Olofsson"
$arr[4,49] = "Roger Olofsson"
$arr[4,50] = $null

$arr[5,0] = 111813872
$arr[5,1] = 56398
$arr[5,2] = "Ovaliderad"
$arr[5,3] = "NT"
$arr[5,4] = 100109
$arr[5,5] = "Tretåig hackspett"
$arr[5,6] = "Picoides tridactylus"
$arr[5,7] = "(Linnaeus, 1758)"
$arr[5,8] = $null
$arr[5,9] = $null
$arr[5,10] = $null
$arr[5,11] = $null
$arr[5,12] = "färska spår"
$arr[5,13] = $null
$arr[5,14] = $null
$arr[5,15] = "Matsdal, granskog, Ås lm"
$arr[5,16] = 540557.5018987871
$arr[5,17] = 7247552.715308581
$arr[5,18] = 10
$arr[5,19] = "Västerbotten"
$arr[5,20] = "Vilhelmina"
$arr[5,21] = "Åsele lappmark"
$arr[5,22] = "Vilhelmina"
$arr[5,23] = $null
$arr[5,24] = "2023-08-18"
$arr[5,25] = "11:00"
$arr[5,26] = "2023-08-18"
$arr[5,27] = "11:00"
$arr[5,28] = $null
$arr[5,29] = $false
$arr[5,30] = $false
$arr[5,31] = $null
$arr[5,32] = $false
$arr[5,33] = $null
$arr[5,34] = $null
$arr[5,35] = $null
$arr[5,36] = $null
$arr[5,37] = $null
$arr[5,38] = $null
$arr[5,39] = $null
$arr[5,40] = $null
$arr[5,41] = $null
$arr[5,42] = $null
$arr[5,43] = $null
$arr[5,44] = $null
$arr[5,45] = $null
$arr[5,46] = $null
$arr[5,47] = $null
$arr[5,48] = "Roger Olofsson"
$arr[5,49] = "Roger Olofsson"
$arr[5,50] = $null

$arr[6,0] = 111814119
$arr[6,1] = 89423
$arr[6,2] = "Ovaliderad"
$arr[6,3] = "NT"
$arr[6,4] = 5432
$arr[6,5] = "Granticka"
$arr[6,6] = "Porodaedalea chrysoloma"
$arr[6,7] = "(Fr.) Fiasson & Niemelä"
$arr[6,8] = $null
$arr[6,9] = $null
$arr[6,10] = $null
$arr[6,11] = $null
$arr[6,12] = $null
$arr[6,13] = $null
$arr[6,14] = $null
$arr[6,15] = "Matsdal, granskog, Ås lm"
$arr[6,16] = 540683.0369185829
$arr[6,17] = 7247576.171207689
$arr[6,18] = 10
$arr[6,19] = "Västerbotten"
$arr[6,20] = "Vilhelmina"
$arr[6,21] = "Åsele lappmark"
$arr[6,22] = "Vilhelmina"
$arr[6,23] = $null
$arr[6,24] = "2023-08-18"
$arr[6,25] = "11:00"
$arr[6,26] = "2023-08-18"
$arr[6,27] = "11:00"
$arr[6,28] = $null
$arr[6,29] = $false
$arr[6,30] = $false
$arr[6,31] = $null
$arr[6,32] = $false
$arr[6,33] = $null
$arr[6,34] = $null
$arr[6,35] = $null
$arr[6,36] = $null
$arr[6,37] = $null
$arr[6,38] = $null
$arr[6,39] = $null
$arr[6,40] = $null
$arr[6,41] = $null
$arr[6,42] = $null
$arr[6,43] = $null
$arr[6,44] = $null
$arr[6,45] = $null
$arr[6,46] = $null
$arr[6,47] = $null
$arr[6,48] = "Roger Olofsson"
$arr[6,49] = "Roger Olofsson"
$arr[6,50] = $null

$arr[7,0] = 111814135
$arr[7,1] = 90087
$arr[7,2] = "Ovaliderad"
$arr[7,3] = "LC"
$arr[7,4] = 3298
$arr[7,5] = "Trådticka"
$arr[7,6] = "Climacocystis borealis"
$arr[7,7] = "(Fr.) Kotl. & Pouzar"
$arr[7,8] = $null
$arr[7,9] = $null
$arr[7,10] = $null
$arr[7,11] = $null
$arr[7,12] = $null
$arr[7,13] = $null
$arr[7,14] = $null
$arr[7,15] = "Matsdal, granskog, Ås lm"
$arr[7,16] = 540661.0419420782
$arr[7,17] = 7247564.172119373
$arr[7,18] = 10
$arr[7,19] = "Västerbotten"
$arr[7,20] = "Vilhelmina"
$arr[7,21] = "Åsele lappmark"
$arr[7,22] = "Vilhelmina"
$arr[7,23] = $null
$arr[7,24] = "2023-08-18"
$arr[7,25] = "11:00"
$arr[7,26] = "2023-08-18"
$arr[7,27] = "11:00"
$arr[7,28] = $null
$arr[7,29] = $false
$arr[7,30] = $false
$arr[7,31] = $null
$arr[7,32] = $false
$arr[7,33] = $null
$arr[7,34] = $null
$arr[7,35] = $null
$arr[7,36] = $null
$arr[7,37] = $null
$arr[7,38] = $null
$arr[7,39] = $null
$arr[7,40] = $null
$arr[7,41] = $null
$arr[7,42] = $null
$arr[7,43] = $null
$arr[7,44] = $null
$arr[7,45] = $null
$arr[7,46] = $null
$arr[7,47] = $null
$arr[7,48] = "Roger Olofsson"
$arr[7,49] = "Roger Olofsson"
$arr[7,50] = $null

$arr[8,0] = 111813938
$arr[8,1] = 89423
$arr[8,2] = "Ovaliderad"
$arr[8,3] = "NT"
$arr[8,4] = 5432
$arr[8,5] = "Granticka"
$arr[8,6] = "Porodaedalea chrysoloma"
$arr[8,7] = "(Fr.) Fiasson & Niemelä"
$arr[8,8] = $null
$arr[8,9] = $null
$arr[8,10] = $null
$arr[8,11] = $null
$arr[8,12] = $null
$arr[8,13] = $null
$arr[8,14] = $null
$arr[8,15] = "Matsdal, granskog, Ås lm"
$arr[8,16] = 540654.849203686
$arr[8,17] = 7247498.096959669
$arr[8,18] = 10
$arr[8,19] = "Västerbotten"
$arr[8,20] = "Vilhelmina"
$arr[8,21] = "Åsele lappmark"
$arr[8,22] = "Vilhelmina"
$arr[8,23] = $null
$arr[8,24] = "2023-08-18"
$arr[8,25] = "11:00"
$arr[8,26] = "2023-08-18"
$arr[8,27] = "11:00"
$arr[8,28] = $null
$arr[8,29] = $false
$arr[8,30] = $false
$arr[8,31] = $null
$arr[8,32] = $false
$arr[8,33] = $null
$arr[8,34] = $null
$arr[8,35] = $null
$arr[8,36] = $null
$arr[8,37] = $null
$arr[8,38] = $null
$arr[8,39] = $null
$arr[8,40] = $null
$arr[8,41] = $null
$arr[8,42] = $null
$arr[8,43] = $null
$arr[8,44] = $null
$arr[8,45] = $null
$arr[8,46] = $null
$arr[8,47] = $null
$arr[8,48] = "Roger Olofsson"
$arr[8,49] = "Roger Olofsson"
$arr[8,50] = $null

$arr[9,0] = 111813975
$arr[9,1] = 89423
$arr[9,2] = "Ovaliderad"
$arr[9,3] = "NT"
$arr[9,4] = 5432
$arr[9,5] = "Granticka"
$arr[9,6] = "Porodaedalea chrysoloma"
$arr[9,7] = "(Fr.) Fiasson & Niemelä"
$arr[9,8] = $null
$arr[9,9] = $null
$arr[9,10] = $null
$arr[9,11] = $null
$arr[9,12] = $null
$arr[9,13] = $null
$arr[9,14] = $null
$arr[9,15] = "Matsdal, granskog, Ås lm"
$arr[9,16] = 540643.7191088985
$arr[9,17] = 7247516.737328541
$arr[9,18] = 10
$arr[9,19] = "Västerbotten"
$arr[9,20] = "Vilhelmina"
$arr[9,21] = "Åsele lappmark"
$arr[9,22] = "Vilhelmina"
$arr[9,23] = $null
$arr[9,24] = "2023-08-18"
$arr[9,25] = "11:00"
$arr[9,26] = "2023-08-18"
$arr[9,27] = "11:00"
$arr[9,28] = $null
$arr[9,29] = $false
$arr[9,30] = $false
$arr[9,31] = $null
$arr[9,32] = $false
$arr[9,33] = $null
$arr[9,34] = $null
$arr[9,35] = $null
$arr[9,36] = $null
$arr[9,37] = $null
$arr[9,38] = $null
$arr[9,39] = $null
$arr[9,40] = $null
$arr[9,41] = $null
$arr[9,42] = $null
$arr[9,43] = $null
$arr[9,44] = $null
$arr[9,45] = $null
$arr[9,46] = $null
$arr[9,47] = $null
$arr[9,48] = "Roger Olofsson"
$arr[9,49] = "Roger Olofsson"
$arr[9,50] = $null

$arr[10,0] = 111813785
$arr[10,1] = 89405
$arr[10,2] = "Ovaliderad"
$arr[10,3] = "NT"
$arr[10,4] = 1202
$arr[10,5] = "Ullticka"
$arr[10,6] = "Phellinidium ferrugineofuscum"
$arr[10,7] = "(P.Karst.) Fiasson & Niemelä"
$arr[10,8] = $null
$arr[10,9] = $null
$arr[10,10] = $null
$arr[10,11] = $null
$arr[10,12] = $null
$arr[10,13] = $null
$arr[10,14] = $null
$arr[10,15] = "Matsdal, granskog, Ås lm"
$arr[10,16] = 540570.9514120822
$arr[10,17] = 7247577.960198429
$arr[10,18] = 10
$arr[10,19] = "Västerbotten"
$arr[10,20] = "Vilhelmina"
$arr[10,21] = "Åsele lappmark"
$arr[10,22] = "Vilhelmina"
$arr[10,23] = $null
$arr[10,24] = "2023-08-18"
$arr[10,25] = "11:00"
$arr[10,26] = "2023-08-18"
$arr[10,27] = "11:00"
$arr[10,28] = $null
$arr[10,29] = $false
$arr[10,30] = $false
$arr[10,31] = $null
$arr[10,32] = $false
$arr[10,33] = $null
$arr[10,34] = $null
$arr[10,35] = $null
$arr[10,36] = $null
$arr[10,37] = $null
$arr[10,38] = $null
$arr[10,39] = $null
$arr[10,40] = $null
$arr[10,41] = $null
$arr[10,42] = $null
$arr[10,43] = $null
$arr[10,44] = $null
$arr[10,45] = $null
$arr[10,46] = $null
$arr[10,47] = $null
$arr[10,48] = "Roger Olofsson"
$arr[10,49] = "Roger Olofsson"
$arr[10,50] = $null

$arr[11,0] = 111813707
$arr[11,1] = 56398
$arr[11,2] = "Ovaliderad"
$arr[11,3] = "NT"
$arr[11,4] = 100109
$arr[11,5] = "Tretåig hackspett"
$arr[11,6] = "Picoides tridactylus"
$arr[11,7] = "(Linnaeus, 1758)"
$arr[11,8] = $null
$arr[11,9] = $null
$arr[11,10] = $null
$arr[11,11] = $null
$arr[11,12] = "färska spår"
$arr[11,13] = $null
$arr[11,14] = $null
$arr[11,15] = "Matsdal, granskog, Ås lm"
$arr[11,16] = 540647.037727406
$arr[11,17] = 7247579.013394679
$arr[11,18] = 10
$arr[11,19] = "Västerbotten"
$arr[11,20] = "Vilhelmina"
$arr[11,21] = "Åsele lappmark"
$arr[11,22] = "Vilhelmina"
$arr[11,23] = $null
$arr[11,24] = "2023-08-18"
$arr[11,25] = "11:00"
$arr[11,26] = "2023-08-18"
$arr[11,27] = "11:00"
$arr[11,28] = $null
$arr[11,29] = $false
$arr[11,30] = $false
$arr[11,31] = $null
$arr[11,32] = $false
$arr[11,33] = $null
$arr[11,34] = $null
$arr[11,35] = "gran"
$arr[11,36] = "Picea abies"
$arr[11,37] = $null
$arr[11,38] = $null
$arr[11,39] = $null
$arr[11,40] = "Picea abies"
$arr[11,41] = $null
$arr[11,42] = $null
$arr[11,43] = $null
$arr[11,44] = $null
$arr[11,45] = $null
$arr[11,46] = $null
$arr[11,47] = $null
$arr[11,48] = "Roger Olofsson"
$arr[11,49] = "Roger Olofsson"
$arr[11,50] = $null

$arr[12,0] = 111825245
$arr[12,1] = 89745
$arr[12,2] = "Ovaliderad"
$arr[12,3] = "VU"
$arr[12,4] = 2062
$arr[12,5] = "Ulltickeporing"
$arr[12,6] = "Skeletocutis brevispora"
$arr[12,7] = "Niemelä"
$arr[12,8] = $null
$arr[12,9] = $null
$arr[12,10] = $null
$arr[12,11] = $null
$arr[12,12] = $null
$arr[12,13] = $null
$arr[12,14] = $null
$arr[12,15] = "Matsdal, granskog, Ås lm"
$arr[12,16] = 540641.3816113799
$arr[12,17] = 7247564.734995849
$arr[12,18] = 10
$arr[12,19] = "Västerbotten"
$arr[12,20] = "Vilhelmina"
$arr[12,21] = "Åsele lappmark"
$arr[12,22] = "Vilhelmina"
$arr[12,23] = $null
$arr[12,24] = "2023-08-13"
$arr[12,25] = "00:00"
$arr[12,26] = "2023-08-13"
$arr[12,27] = "00:00"
$arr[12,28] = "Färskt exemplar. Kollekt tog och torkades, gulnade."
$arr[12,29] = $false
$arr[12,30] = $false
$arr[12,31] = $null
$arr[12,32] = $false
$arr[12,33] = $null
$arr[12,34] = $null
$arr[12,35] = "ullticka"
$arr[12,36] = "Phellinidium ferrugineofuscum"
$arr[12,37] = $null
$arr[12,38] = $null
$arr[12,39] = $null
$arr[12,40] = "Phellinidium ferrugineofuscum"
$arr[12,41] = $null
$arr[12,42] = $null
$arr[12,43] = $null
$arr[12,44] = $null
$arr[12,45] = $null
$arr[12,46] = $null
$arr[12,47] = $null
$arr[12,48] = "Roger Olofsson"
$arr[12,49] = "Roger Olofsson"
$arr[12,50] = $null

$arr[13,0] = 111825158
$arr[13,1] = 89405
$arr[13,2] = "Ovaliderad"
$arr[13,3] = "NT"
$arr[13,4] = 1202
$arr[13,5] = "Ullticka"
$arr[13,6] = "Phellinidium ferrugineofuscum"
$arr[13,7] = "(P.Karst.) Fiasson & Niemelä"
$arr[13,8] = $null
$arr[13,9] = $null
$arr[13,10] = $null
$arr[13,11] = $null
$arr[13,12] = $null
$arr[13,13] = $null
$arr[13,14] = $null
$arr[13,15] = "Matsdal, granskog, Ås lm"
$arr[13,16] = 540641.3816113799
$arr[13,17] = 7247564.734995849
$arr[13,18] = 10
$arr[13,19] = "Västerbotten"
$arr[13,20] = "Vilhelmina"
$arr[13,21] = "Åsele lappmark"
$arr[13,22] = "Vilhelmina"
$arr[13,23] = $null
$arr[13,24] = "2023-08-13"
$arr[13,25] = "00:00"
$arr[13,26] = "2023-08-13"
$arr[13,27] = "00:00"
$arr[13,28] = "Med ulltickeporing"
$arr[13,29] = $false
$arr[13,30] = $false
$arr[13,31] = $null
$arr[13,32] = $false
$arr[13,33] = $null
$arr[13,34] = $null
$arr[13,35] = "gran"
$arr[13,36] = "Picea abies"
$arr[13,37] = $null
$arr[13,38] = $null
$arr[13,39] = $null
$arr[13,40] = "Picea abies"
$arr[13,41] = $null
$arr[13,42] = $null
$arr[13,43] = $null
$arr[13,44] = $null
$arr[13,45] = $null
$arr[13,46] = $null
$arr[13,47] = $null
$arr[13,48] = "Roger Olofsson"
$arr[13,49] = "Roger Olofsson"
$arr[13,50] = $null

$arr[14,0] = 111825098
$arr[14,1] = 89845
$arr[14,2] = "Ovaliderad"
$arr[14,3] = "VU"
$arr[14,4] = 1209
$arr[14,5] = "Rynkskinn"
$arr[14,6] = "Phlebia centrifuga"
$arr[14,7] = "P.Karst."
$arr[14,8] = $null
$arr[14,9] = $null
$arr[14,10] = $null
$arr[14,11] = $null
$arr[14,12] = $null
$arr[14,13] = $null
$arr[14,14] = $null
$arr[14,15] = "Matsdal, granskog, Ås lm"
$arr[14,16] = 540641.3816113799
$arr[14,17] = 7247564.734995849
$arr[14,18] = 10
$arr[14,19] = "Västerbotten"
$arr[14,20] = "Vilhelmina"
$arr[14,21] = "Åsele lappmark"
$arr[14,22] = "Vilhelmina"
$arr[14,23] = $null
$arr[14,24] = "2023-08-13"
$arr[14,25] = "00:00"
$arr[14,26] = "2023-08-13"
$arr[14,27] = "00:00"
$arr[14,28] = "Färskt. På granlåga med minst 25 rosentickor, ullticka och ulltickeporing.  Någon gul slemsvamp? på rynkskinnet"
$arr[14,29] = $false
$arr[14,30] = $false
$arr[14,31] = $null
$arr[14,32] = $false
$arr[14,33] = $null
$arr[14,34] = $null
$arr[14,35] = "gran"
$arr[14,36] = "Picea abies"
$arr[14,37] = $null
$arr[14,38] = $null
$arr[14,39] = $null
$arr[14,40] = "Picea abies"
$arr[14,41] = $null
$arr[14,42] = $null
$arr[14,43] = $null
$arr[14,44] = $null
$arr[14,45] = $null
$arr[14,46] = $null
$arr[14,47] = $null
$arr[14,48] = "Roger Olofsson"
$arr[14,49] = "Roger Olofsson"
$arr[14,50] = $null

$arr[15,0] = 111825340
$arr[15,1] = 89686
$arr[15,2] = "Ovaliderad"
$arr[15,3] = "NT"
$arr[15,4] = 658
$arr[15,5] = "Rosenticka"
$arr[15,6] = "Rhodofomes roseus"
$arr[15,7] = "(Alb. & Schwein.) Kotl. & Pouzar"
$arr[15,8] = "25"
$arr[15,9] = $null
$arr[15,10] = $null
$arr[15,11] = $null
$arr[15,12] = $null
$arr[15,13] = $null
$arr[15,14] = $null
$arr[15,15] = "Matsdal, granskog, Ås lm"
$arr[15,16] = 540641.3816113799
$arr[15,17] = 7247564.734995849
$arr[15,18] = 10
$arr[15,19] = "Västerbotten"
$arr[15,20] = "Vilhelmina"
$arr[15,21] = "Åsele lappmark"
$arr[15,22] = "Vilhelmina"
$arr[15,23] = $null
$arr[15,24] = "2023-08-13"
$arr[15,25] = "00:00"
$arr[15,26] = "2023-08-13"
$arr[15,27] = "00:00"
$arr[15,28] = "Minst 25 rosentickor på denna låga. Ullticka, ulltickeporing, rynkskinn på samma låga"
$arr[15,29] = $false
$arr[15,30] = $false
$arr[15,31] = $null
$arr[15,32] = $false
$arr[15,33] = $null
$arr[15,34] = $null
$arr[15,35] = $null
$arr[15,36] = $null
$arr[15,37] = $null
$arr[15,38] = $null
$arr[15,39] = $null
$arr[15,40] = $null
$arr[15,41] = $null
$arr[15,42] = $null
$arr[15,43] = $null
$arr[15,44] = $null
$arr[15,45] = $null
$arr[15,46] = $null
$arr[15,47] = $null
$arr[15,48] = "Roger Olofsson"
$arr[15,49] = "Roger Olofsson"
$arr[15,50] = $null

$ws.Range("A2:AY17").Value = $arr

# Re-apply the cells whose text would otherwise be auto-coerced into
# a number/date/time by Excel, forcing them to stay as literal text,
# then reset the style back to Normal so no formatting diff remains.
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2023-08-18"
$ws.Range("Y2").Style = "Normal"
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "11:00"
$ws.Range("Z2").Style = "Normal"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2023-08-18"
$ws.Range("AA2").Style = "Normal"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "11:00"
$ws.Range("AB2").Style = "Normal"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-08-18"
$ws.Range("Y3").Style = "Normal"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "11:00"
$ws.Range("Z3").Style = "Normal"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-08-18"
$ws.Range("AA3").Style = "Normal"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "11:00"
$ws.Range("AB3").Style = "Normal"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-08-18"
$ws.Range("Y4").Style = "Normal"
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "11:00"
$ws.Range("Z4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-08-18"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "11:00"
$ws.Range("AB4").Style = "Normal"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-08-18"
$ws.Range("Y5").Style = "Normal"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "11:00"
$ws.Range("Z5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-08-18"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "11:00"
$ws.Range("AB5").Style = "Normal"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-08-18"
$ws.Range("Y6").Style = "Normal"
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "11:00"
$ws.Range("Z6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-08-18"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "11:00"
$ws.Range("AB6").Style = "Normal"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-08-18"
$ws.Range("Y7").Style = "Normal"
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = "11:00"
$ws.Range("Z7").Style = "Normal"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-08-18"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AB7").NumberFormat = "@"
$ws.Range("AB7").Value = "11:00"
$ws.Range("AB7").Style = "Normal"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2023-08-18"
$ws.Range("Y8").Style = "Normal"
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = "11:00"
$ws.Range("Z8").Style = "Normal"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2023-08-18"
$ws.Range("AA8").Style = "Normal"
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = "11:00"
$ws.Range("AB8").Style = "Normal"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-08-18"
$ws.Range("Y9").Style = "Normal"
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = "11:00"
$ws.Range("Z9").Style = "Normal"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2023-08-18"
$ws.Range("AA9").Style = "Normal"
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = "11:00"
$ws.Range("AB9").Style = "Normal"
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2023-08-18"
$ws.Range("Y10").Style = "Normal"
$ws.Range("Z10").NumberFormat = "@"
$ws.Range("Z10").Value = "11:00"
$ws.Range("Z10").Style = "Normal"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "2023-08-18"
$ws.Range("AA10").Style = "Normal"
$ws.Range("AB10").NumberFormat = "@"
$ws.Range("AB10").Value = "11:00"
$ws.Range("AB10").Style = "Normal"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2023-08-18"
$ws.Range("Y11").Style = "Normal"
$ws.Range("Z11").NumberFormat = "@"
$ws.Range("Z11").Value = "11:00"
$ws.Range("Z11").Style = "Normal"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "2023-08-18"
$ws.Range("AA11").Style = "Normal"
$ws.Range("AB11").NumberFormat = "@"
$ws.Range("AB11").Value = "11:00"
$ws.Range("AB11").Style = "Normal"
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "2023-08-18"
$ws.Range("Y12").Style = "Normal"
$ws.Range("Z12").NumberFormat = "@"
$ws.Range("Z12").Value = "11:00"
$ws.Range("Z12").Style = "Normal"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "2023-08-18"
$ws.Range("AA12").Style = "Normal"
$ws.Range("AB12").NumberFormat = "@"
$ws.Range("AB12").Value = "11:00"
$ws.Range("AB12").Style = "Normal"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2023-08-18"
$ws.Range("Y13").Style = "Normal"
$ws.Range("Z13").NumberFormat = "@"
$ws.Range("Z13").Value = "11:00"
$ws.Range("Z13").Style = "Normal"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2023-08-18"
$ws.Range("AA13").Style = "Normal"
$ws.Range("AB13").NumberFormat = "@"
$ws.Range("AB13").Value = "11:00"
$ws.Range("AB13").Style = "Normal"
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2023-08-13"
$ws.Range("Y14").Style = "Normal"
$ws.Range("Z14").NumberFormat = "@"
$ws.Range("Z14").Value = "00:00"
$ws.Range("Z14").Style = "Normal"
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = "2023-08-13"
$ws.Range("AA14").Style = "Normal"
$ws.Range("AB14").NumberFormat = "@"
$ws.Range("AB14").Value = "00:00"
$ws.Range("AB14").Style = "Normal"
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2023-08-13"
$ws.Range("Y15").Style = "Normal"
$ws.Range("Z15").NumberFormat = "@"
$ws.Range("Z15").Value = "00:00"
$ws.Range("Z15").Style = "Normal"
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "2023-08-13"
$ws.Range("AA15").Style = "Normal"
$ws.Range("AB15").NumberFormat = "@"
$ws.Range("AB15").Value = "00:00"
$ws.Range("AB15").Style = "Normal"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2023-08-13"
$ws.Range("Y16").Style = "Normal"
$ws.Range("Z16").NumberFormat = "@"
$ws.Range("Z16").Value = "00:00"
$ws.Range("Z16").Style = "Normal"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2023-08-13"
$ws.Range("AA16").Style = "Normal"
$ws.Range("AB16").NumberFormat = "@"
$ws.Range("AB16").Value = "00:00"
$ws.Range("AB16").Style = "Normal"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "25"
$ws.Range("I17").Style = "Normal"
$ws.Range("Y17").NumberFormat = "@"
$ws.Range("Y17").Value = "2023-08-13"
$ws.Range("Y17").Style = "Normal"
$ws.Range("Z17").NumberFormat = "@"
$ws.Range("Z17").Value = "00:00"
$ws.Range("Z17").Style = "Normal"
$ws.Range("AA17").NumberFormat = "@"
$ws.Range("AA17").Value = "2023-08-13"
$ws.Range("AA17").Style = "Normal"
$ws.Range("AB17").NumberFormat = "@"
$ws.Range("AB17").Value = "00:00"
$ws.Range("AB17").Style = "Normal"

Write-Output "rows written"
